$wb = $excel.ActiveWorkbook

# --- Awards sheet: insert a new "EMIP Cover Showcase Winner" row for 2020 ---
$awards = $wb.Worksheets.Item("awards")
$awards.Rows.Item(2).Insert()
$awards.Range("A2").Value = "Educational Measurement: Issues and Practice Cover Showcase Winner"
$awards.Range("B2").Value = "National Council on Measurement in Education"
$awards.Range("C2").Value = 2020
$awards.Range("E2").Value = $true

# --- Grants sheet: bump the grant year ranges ---
$grants = $wb.Worksheets.Item("grants")
$grants.Range("E2").Value = 2021
$grants.Range("F2").Value = 2023
$grants.Range("E3").Value = 2021
$grants.Range("F3").Value = 2025
$grants.Range("E5").Value = 2020
$grants.Range("F5").Value = 2022

# --- Selection / active-sheet bookkeeping to mirror the authored workbook state ---
$education = $wb.Worksheets.Item("education")
$education.Range("I9").Select() | Out-Null

$grants.Range("B4").Select() | Out-Null

$awards.Range("C3").Select() | Out-Null
$awards.Activate()
